# Scheduled market-data refresh: update currentAveragePrice*/LevePrice*/LeveProfit*
# columns (H:N) for the leves whose linked item prices moved, across all crafting
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW). WVR is unaffected this run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 583
$ws.Range("I12").Value = 512.4286
$ws.Range("J12").Value = 706.5
$ws.Range("K12").Value = 512.4286
$ws.Range("L12").Value = 706.5
$ws.Range("M12").Value = -342.4286
$ws.Range("N12").Value = -1046.5

$ws.Range("H19").Value = 966.4167
$ws.Range("I19").Value = 1235
$ws.Range("J19").Value = 832.125
$ws.Range("K19").Value = 1235
$ws.Range("L19").Value = 832.125
$ws.Range("M19").Value = -1060
$ws.Range("N19").Value = -1182.125

$ws.Range("H88").Value = 1367
$ws.Range("I88").Value = 1419.5
$ws.Range("J88").Value = 1304
$ws.Range("K88").Value = 1419.5
$ws.Range("L88").Value = 1304
$ws.Range("M88").Value = -1013.5
$ws.Range("N88").Value = -2116

$ws.Range("H91").Value = 1367
$ws.Range("I91").Value = 1419.5
$ws.Range("J91").Value = 1304
$ws.Range("K91").Value = 1419.5
$ws.Range("L91").Value = 1304
$ws.Range("M91").Value = -15.5
$ws.Range("N91").Value = -4112

$ws.Range("H112").Value = 2380.4443
$ws.Range("I112").Value = 3163.3333
$ws.Range("J112").Value = 1989
$ws.Range("K112").Value = 9489.999899999999
$ws.Range("L112").Value = 5967
$ws.Range("M112").Value = -8381.999899999999
$ws.Range("N112").Value = -8183

$ws.Range("H116").Value = 2910.5
$ws.Range("I116").Value = 2788.25
$ws.Range("J116").Value = 3399.5
$ws.Range("K116").Value = 2788.25
$ws.Range("L116").Value = 3399.5
$ws.Range("M116").Value = 653.75
$ws.Range("N116").Value = -10283.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 999.4666999999999
$ws.Range("I2").Value = 1028.4286
$ws.Range("J2").Value = 594
$ws.Range("K2").Value = 1028.4286
$ws.Range("L2").Value = 594
$ws.Range("M2").Value = -915.4286
$ws.Range("N2").Value = -820

$ws.Range("H32").Value = 836.375
$ws.Range("I32").Value = 836.375
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 836.375
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -549.375

$ws.Range("H116").Value = 999.4666999999999
$ws.Range("I116").Value = 1028.4286
$ws.Range("J116").Value = 594
$ws.Range("K116").Value = 1028.4286
$ws.Range("L116").Value = 594
$ws.Range("M116").Value = 1265.5714
$ws.Range("N116").Value = -5182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 999.4666999999999
$ws.Range("I3").Value = 1028.4286
$ws.Range("J3").Value = 594
$ws.Range("K3").Value = 1028.4286
$ws.Range("L3").Value = 594
$ws.Range("M3").Value = -914.4286
$ws.Range("N3").Value = -822

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 15002490
$ws.Range("I4").Value = 15002490
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 15002490
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -15002378

$ws.Range("H31").Value = 8658.9
$ws.Range("I31").Value = 4538
$ws.Range("J31").Value = 10032.533
$ws.Range("K31").Value = 4538
$ws.Range("L31").Value = 10032.533
$ws.Range("M31").Value = -4243
$ws.Range("N31").Value = -10622.533

$ws.Range("H34").Value = 8658.9
$ws.Range("I34").Value = 4538
$ws.Range("J34").Value = 10032.533
$ws.Range("K34").Value = 4538
$ws.Range("L34").Value = 10032.533
$ws.Range("M34").Value = -4336
$ws.Range("N34").Value = -10436.533

$ws.Range("H99").Value = 1924.6923
$ws.Range("I99").Value = 1853.5
$ws.Range("J99").Value = 1985.7142
$ws.Range("K99").Value = 1853.5
$ws.Range("L99").Value = 1985.7142
$ws.Range("M99").Value = -355.5
$ws.Range("N99").Value = -4981.7142

$ws.Range("H126").Value = 1924.6923
$ws.Range("I126").Value = 1853.5
$ws.Range("J126").Value = 1985.7142
$ws.Range("K126").Value = 5560.5
$ws.Range("L126").Value = 5957.142599999999
$ws.Range("M126").Value = -3090.5
$ws.Range("N126").Value = -10897.1426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 878.6
$ws.Range("I5").Value = 777
$ws.Range("J5").Value = 1031
$ws.Range("K5").Value = 2331
$ws.Range("L5").Value = 3093
$ws.Range("M5").Value = -2219
$ws.Range("N5").Value = -3317

$ws.Range("H19").Value = 4999
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 4999
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 14997
$ws.Range("N19").Value = -15345

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H68").Value = 2541.4
$ws.Range("I68").Value = 1475
$ws.Range("J68").Value = 3252.3333
$ws.Range("K68").Value = 4425
$ws.Range("L68").Value = 9756.999899999999
$ws.Range("M68").Value = -3614
$ws.Range("N68").Value = -11378.9999

$ws.Range("H71").Value = 2541.4
$ws.Range("I71").Value = 1475
$ws.Range("J71").Value = 3252.3333
$ws.Range("K71").Value = 13275
$ws.Range("L71").Value = 29270.9997
$ws.Range("M71").Value = -9219
$ws.Range("N71").Value = -37382.9997

$ws.Range("H80").Value = 4346.364
$ws.Range("I80").Value = 4058.7646
$ws.Range("J80").Value = 5324.2
$ws.Range("K80").Value = 12176.2938
$ws.Range("L80").Value = 15972.6
$ws.Range("M80").Value = -11240.2938
$ws.Range("N80").Value = -17844.6

$ws.Range("H83").Value = 4346.364
$ws.Range("I83").Value = 4058.7646
$ws.Range("J83").Value = 5324.2
$ws.Range("K83").Value = 36528.8814
$ws.Range("L83").Value = 47917.8
$ws.Range("M83").Value = -31848.8814
$ws.Range("N83").Value = -57277.8

$ws.Range("H113").Value = 1423.2
$ws.Range("I113").Value = 479
$ws.Range("J113").Value = 1659.25
$ws.Range("K113").Value = 1437
$ws.Range("L113").Value = 4977.75
$ws.Range("M113").Value = 733
$ws.Range("N113").Value = -9317.75

$ws.Range("H129").Value = 840.5714
$ws.Range("I129").Value = 874.75
$ws.Range("J129").Value = 795
$ws.Range("K129").Value = 2624.25
$ws.Range("L129").Value = 2385
$ws.Range("M129").Value = 2375.75
$ws.Range("N129").Value = -12385

$ws.Range("H135").Value = 878.6
$ws.Range("I135").Value = 777
$ws.Range("J135").Value = 1031
$ws.Range("K135").Value = 6993
$ws.Range("L135").Value = 9279
$ws.Range("M135").Value = -4458
$ws.Range("N135").Value = -14349

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 210.6923
$ws.Range("I2").Value = 155.33333
$ws.Range("J2").Value = 258.14285
$ws.Range("K2").Value = 155.33333
$ws.Range("L2").Value = 258.14285
$ws.Range("M2").Value = -42.33332999999999
$ws.Range("N2").Value = -484.14285

$ws.Range("H46").Value = 14222.223
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 17428.572
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 17428.572
$ws.Range("M46").Value = -2844
$ws.Range("N46").Value = -17740.572

$ws.Range("H57").Value = 6513.75
$ws.Range("I57").Value = 3685
$ws.Range("J57").Value = 15000
$ws.Range("K57").Value = 3685
$ws.Range("L57").Value = 15000
$ws.Range("M57").Value = -2865
$ws.Range("N57").Value = -16640

$ws.Range("H70").Value = 2000
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -1730

$ws.Range("H73").Value = 2000
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -1064

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1550
$ws.Range("I22").Value = 2100
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 2100
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -1805
$ws.Range("N22").Value = -1590

$ws.Range("H27").Value = 1550
$ws.Range("I27").Value = 2100
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 2100
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -1993
$ws.Range("N27").Value = -1214

$ws.Range("H46").Value = 5909.091
$ws.Range("I46").Value = 6000
$ws.Range("J46").Value = 5857.143
$ws.Range("K46").Value = 6000
$ws.Range("L46").Value = 5857.143
$ws.Range("M46").Value = -5812
$ws.Range("N46").Value = -6233.143
